{"js": "// Add the \"O.K\" note into the END column of Ehab Kadhum's row (the\n// first data row of the status table), which was previously empty.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Row 0 is the header (NO | NAME | TASK | END); row 1 is Ehab Kadhum's row.\nconst targetRow = rows.items[1];\nconst cells = targetRow.cells;\ncells.load(\"items\");\nawait context.sync();\n\n// Column 3 (0-based) is the \"END\" column.\nconst endCell = cells.items[3];\nendCell.body.insertText(\"O.K\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Add the \"O.K\" note into the END column of Ehab Kadhum's row (the\n# first data row of the status table), which was previously empty.\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Row 1 is the header (NO | NAME | TASK | END); row 2 is Ehab Kadhum's row.\n# Column 4 is the \"END\" column.\n$cell = $table.Cell(2, 4)\n$cell.Range.Text = \"O.K\"\n"}
